# Uppercase the terminal names in column K (rows 2-58) of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 11).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 11)  # column K = 11
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val.ToString().ToUpper()
    }
}
